# The commit swaps the contents of ppt/theme/theme1.xml (the slide-master's
# "Integral" theme) and ppt/theme/theme2.xml (the notes-master's default
# "Office Theme"), so the slide master ends up using the plain Office theme
# colors (and vice-versa for the notes master).
#
# In this deck the two themes already share an identical <a:fontScheme> and
# <a:fmtScheme>; the only real content difference between them is the
# <a:clrScheme> (the 12 theme colors). We reproduce that swap through the
# PowerPoint object model by writing the "Office Theme" palette onto the
# slide master's theme color scheme.

$p  = $ppt.ActivePresentation
$sm = $p.SlideMaster
$cs = $sm.Theme.ThemeColorScheme

# Target palette (the stock "Office Theme" colors), in the
# MsoThemeColorSchemeIndex order used by ThemeColorScheme.Item:
#  1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1, 6 accent2, 7 accent3,
#  8 accent4, 9 accent5, 10 accent6, 11 hlink, 12 folHlink
$officeColors = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

for ($i = 1; $i -le 12; $i++) {
    $hex = $officeColors[$i - 1]
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    $cs.Item($i).RGB = $r + ($g * 256) + ($b * 65536)
}

Write-Host "Slide master theme colors updated to the Office Theme palette."
